$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.749.58'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '2.038.75'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''227.55'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = '''0.608'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("D7").Value = '''60.19'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("D10").Value = '''0.0828'
$ws.Range("E10").Value = '  +2.28%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '2.337.87'
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("D13").Value = '''14.59'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '''21.16'
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '''0.777'
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("D16").Value = '''5.35'
$ws.Range("E16").Value = '  +2.59%  '
$ws.Range("D17").Value = '2.033.92'
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = '37.681.61'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("D20").Value = '''69.46'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '''223.98'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '''2.42'
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("E25").Value = '  +4.23%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '''9.38'
$ws.Range("E26").Value = '  +3.01%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '''167.78'
$ws.Range("E27").Value = '  +1.28%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '''18.79'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").Value = '''1.27'
$ws.Range("E30").Value = '  -1.35%  '
$ws.Range("D31").Value = '''0.121'
$ws.Range("E31").Value = '  +1.04%  '
$ws.Range("E32").Value = '  +8.57%  '
$ws.Range("D33").Value = '''4.38'
$ws.Range("E33").Value = '  -1.16%  '
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("D35").Value = '''4.48'
$ws.Range("E35").Value = '  -0.26%  '
$ws.Range("E36").Value = '  +4.51%  '
$ws.Range("E37").Value = '  +4.14%  '
$ws.Range("D38").Value = '''3.42'
$ws.Range("E38").Value = '  +6.74%  '
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("D40").Value = '''18.07'
$ws.Range("E40").Value = '  +10.05%  '
$ws.Range("D41").Value = '1.528.13'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '''97.31'
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("D45").Value = '''0.0910'
$ws.Range("E45").Value = '  -0.59%  '
$ws.Range("D46").Value = '''4.15'
$ws.Range("E46").Value = '  +4.13%  '
$ws.Range("D47").Value = '''1.10'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").Value = '''7.06'
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").Value = '2.226.28'
$ws.Range("E51").Value = '  +0.60%  '
